# Apply the diff:
#  - Slide 5 "Content Placeholder 2": resize/move (off/ext) only.
#  - Slide 6 "Content Placeholder 2": resize/move, shrink body text from
#    14pt to 12pt, and append a new level-2 bullet paragraph.

$p = $ppt.ActivePresentation

# --- Slide 5 ---------------------------------------------------------
$s5  = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(3)          # "Content Placeholder 2"

# Target EMU: off x=533400 y=3193752, ext cx=7924800 cy=1327365
# (point values nudged by a hair above the exact EMU/12700 quotient so
#  the host's internal f32 storage still truncates back to the right EMU)
$sh5.Left   = 42.0
$sh5.Top    = 251.47653963307087
$sh5.Width  = 624.0
$sh5.Height = 104.51692913385827

# --- Slide 6 ---------------------------------------------------------
$s6  = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(3)          # "Content Placeholder 2"

# Target EMU: off x=381000 y=2742429, ext cx=8382000 cy=2123658
$sh6.Left   = 30.0
$sh6.Top    = 215.93929293858267
$sh6.Width  = 660.0
$sh6.Height = 167.2171707543307

$tr = $sh6.TextFrame.TextRange

# Shrink every existing paragraph's run text from 14pt to 12pt.
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $tr.Paragraphs($i).Font.Size = 12
}

# Append a new level-2 bullet after the last paragraph, inheriting the
# same paragraph formatting (lnSpc/spcBef/indent level) as paragraph 5.
$newPara = $tr.InsertAfter([char]13 + "Loopback probe packets reach the intended reflector node")
$newPara.Font.Size = 12
